$d = $word.ActiveDocument

# 1. Insert a new run containing "ex" right before the existing "Overview"
#    run in the very first paragraph (heading), so the visible text becomes
#    "exOverview".
$d.Content.Find.Execute("Overview", $false, $false, $false, $false, $false,
                         $true, 1, $false, "exOverview", 2)

# 2. Tighten up every table's left indent / left cell padding:
#    tblInd -49 -> -59 dxa  (i.e. -2.45pt -> -2.95pt)
#    tblCellMar/left and tcMar/left 40 -> 30 dxa (i.e. 2pt -> 1.5pt)
foreach ($t in $d.Tables) {
    $t.Rows.LeftIndent = -2.95
    $t.LeftPadding = 1.5
    $t.Cell(1, 1).LeftPadding = 1.5
}
